$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3-6 (columns A,B,E,F,G,H,Q,R,AC) are cyclically rotated:
# new row3 <- old row6, new row4 <- old row3, new row5 <- old row4, new row6 <- old row5.
# Capture old values first, then write them back in rotated order so we don't
# clobber data we still need to read.

$cols = @("A","B","E","F","G","H","Q","R","AC")

$old = @{}
foreach ($r in 3..6) {
    $old[$r] = @{}
    foreach ($c in $cols) {
        $old[$r][$c] = $ws.Range("$c$r").Value()
    }
}

# mapping: new row -> source old row
$map = @{ 3 = 6; 4 = 3; 5 = 4; 6 = 5 }

foreach ($newr in 3..6) {
    $srcr = $map[$newr]
    foreach ($c in $cols) {
        $ws.Range("$c$newr").Value = $old[$srcr][$c]
    }
}
